$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (Keyword, Correlation)
$updates = @{
    2  = @("invest",    -0.6263)
    4  = @("uncertain",  0.08019999999999999)
    7  = @("inflation", -0.1311)
    8  = @("trade",      0.1349)
    9  = @("interest",  -0.4211)
    10 = @("uncertain", -0.6423)
    11 = @("invest",    -0.0548)
    12 = @("trade",      0.8954)
    13 = @("uncertain", -0.1748)
    14 = @("interest",   0.4966)
    15 = @("invest",     0.4175)
    17 = @("trade",     -0.1808)
    18 = @("inflation",  0.5832000000000001)
    19 = @("interest",   0.8001)
    20 = @("invest",    -0.1405)
    21 = @("uncertain", -0.2575)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
}
